$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New 1,000-row dataset refresh (custom accuracy pass) ---
# Replace the 4 data rows (rows 2-5) with the refreshed readings and
# drop the now-unused 6th sample row.

$rowsData = @(
  @(45084.50694444445,1.29,0.735,0.124,4.827,1.423,0.979,2.643,1.329,0.541,1.443,2.279,1.076,0.423,1.072,3.872,0.514,0,0.001,14.386,3.712,2.121,1.552,1.645,1.823,2.314,0.973,0.721,1.589,1.435,2.091,2.771,0.187,2.45),
  @(45084.51388888889,18.888,14.094,0.716,41.961,33.589,14.968,50.407,22.863,10.29,15.649,17.139,17.456,4.822,14.776,22.27,12.033,0.176,0.354,221.15,41.832,14.228,28,15.2,2.702,26.37,12.209,10.614,12.906,17.582,0.8169999999999999,45.952,7.592,17.701),
  @(45084.52083333334,13.712,10.241,0.537,30.378,24.412,10.891,43.637,16.617,7.512,11.349,12.436,12.671,3.501,10.715,16.108,8.720000000000001,0.132,0.239,158.386,30.526,10.291,20.456,11.038,1.9,21.517,8.867000000000001,7.698,9.361000000000001,12.766,0.505,40.082,5.521,12.81),
  @(45084.52777777778,23.38,17.52,0.85,51.25,41.86,18.48,70.53,28.32,12.76,19.15,20.77,21.59,5.92,18.28,26.72,15.01,0.21,0.5600000000000001,272.91,51.48,17.17,34.72,18.47,2.76,34.81,15.01,13.09,15.65,21.59,0.37,64.2,9.529999999999999,21.43)
)

for ($r = 0; $r -lt $rowsData.Length; $r++) {
    $rowVals = $rowsData[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}

# The old 6th sample row is no longer part of the refreshed dataset.
$ws.Rows.Item(6).Delete()

# A few columns need to widen by one character to fit the new values.
$ws.Range("C:C").ColumnWidth = 43/6
$ws.Range("G:G").ColumnWidth = 43/6
$ws.Range("K:M").ColumnWidth = 43/6
$ws.Range("O:Q").ColumnWidth = 43/6
$ws.Range("V:V").ColumnWidth = 43/6
$ws.Range("X:X").ColumnWidth = 43/6
$ws.Range("AA:AC").ColumnWidth = 43/6
$ws.Range("AH:AH").ColumnWidth = 43/6
